# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# Rows 183-186 (match records) get their betting-data columns (B, F:AC)
# cyclically rotated while the rank column A stays put:
#   new row 183 <- old row 185
#   new row 184 <- old row 183
#   new row 185 <- old row 186
#   new row 186 <- old row 184

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for columns B and F:AC on the four affected rows
# (column order matches the worksheet layout: B, F, G, H, I, J, K, L, M, N, O,
#  P, Q, R, S, T, U, V, W, X, Y, Z, AA, AB, AC)
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
$rows = @(183,184,185,186)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# target row -> source row (data comes FROM the source row's old/"before" values)
$mapping = @{ 183 = 185; 184 = 183; 185 = 186; 186 = 184 }

foreach ($target in $rows) {
    $source = $mapping[$target]
    $srcData = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $srcData[$c]
    }
}

Write-Host "Rotation applied."
